# PHILIBERT_JournalDeBord.xlsx edit script
# - Extends the "authentification stateless" comment (E9) with a mention of JSON Web Tokens
# - Fixes C20 (durée) from the text "0..5" to the numeric value 0.5
# - Fills in three previously-blank rows (21, 22, 23) with new journal entries:
#     * Analyse de faisabilité (row 21)
#     * Suite/fin de l'analyse concurrentielle (row 22)
#     * Début de la planification détaillée (row 23)
# - Updates the active selection to E24 to reflect the new bottom of the filled data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: extend the stateless-authentication reflection comment ---
$ws.Range('E9').Value = 'C''est la première fois que je rencontre le problème de réalisé une authentification stateless, qui ne repose donc sur aucun stockage côté serveur. Je trouve que certains problèmes sont plus simplement résolus lorsque une tel architecture est mise en place, comme par exemple, l''attribution de scope aux ressources. D''autres problèmes sont cependant présent avec une authetification stateless, comme la mise en place de JSON Web Tokens'
$ws.Rows.Item(9).RowHeight = 120

# --- Row 18: tighten the row height to the actual wrapped-text size ---
$ws.Rows.Item(18).RowHeight = 104.25

# --- Row 20: fix the duration cell, it was accidentally entered as text "0..5" ---
$ws.Range('C20').Value = 0.5
$ws.Rows.Item(20).RowHeight = 74.25

# --- Rows 21-23 (were blank placeholders): carry over row 20's cell formatting first ---
$ws.Range('A20:E20').Copy()
$ws.Range('A21:E23').PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 21: Analyse de faisabilité ---
$ws.Range('A21').Value = $ws.Range('A20').Value()
$ws.Range('B21').Value = 'Écriture de l''analyse de faisabilité, je défini la faisabilité système ainsi que la faisabilité logicielle, aucun problème de faisabilité ne devrait être rencontré '
$ws.Range('C21').Value = 0.5
$ws.Range('D21').Value = $ws.Range('D20').Value()
$ws.Rows.Item(21).RowHeight = 90

# --- Row 22: suite de l'analyse concurrentielle ---
$ws.Range('A22').Value = $ws.Range('A20').Value()
$ws.Range('B22').Value = 'Je termine mon analyse concurrentielle, je termine de remplir mon tableau comparatif et j''ajoute encore quelques points à comparés'
$ws.Range('C22').Value = 0.5
$ws.Range('D22').Value = $ws.Range('D20').Value()
$ws.Rows.Item(22).RowHeight = 75

# --- Row 23: début de la planification détaillée ---
$ws.Range('A23').Value = $ws.Range('A20').Value()
$ws.Range('B23').Value = 'Début de la plannification détaillée'
$ws.Range('C23').Value = 1.5
$ws.Range('D23').Value = $ws.Range('D20').Value()
$ws.Range('E23').Value = 'Je commence la plannification détaillée avec comme modèle la plannification initiale, il y a certains point que je n''arrive pas à détailler plus, vu que ces points sont nouveau pour moi, je n''arrive donc pas à déterminer toutes les étapes nécessaire pour effectuer ces tâches dans la plannifaction à ce moment du projet'
$ws.Rows.Item(23).RowHeight = 90

# --- Update the active selection to reflect scrolling down to the newly filled rows ---
$ws.Range('E24').Select()
